$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.270.62"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "1.664.49"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  +0.88%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.99"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5303"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.011"
$ws.Range("E7").Value = "  +0.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2636"
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06363"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.58"
$ws.Range("E10").Value = "  +0.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07855"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.565"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").Value = "1.665.66"
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("D14").Value = "1.892.55"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5532"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").Value = "0.0₅8176"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.63"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.011"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.672"
$ws.Range("E19").Value = "  +2.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.87"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.23"
$ws.Range("E21").Value = "  +1.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.046"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.012"
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.02"
$ws.Range("E24").Value = "  +2.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1224"
$ws.Range("E25").Value = "  -1.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.214"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.15"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.488"
$ws.Range("E28").Value = "  +3.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05997"
$ws.Range("E29").Value = "  +1.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.281"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.586"
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.289"
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.619"
$ws.Range("E33").Value = "  +3.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9631"
$ws.Range("E34").Value = "  +1.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.829"
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.427"
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5802"
$ws.Range("E37").Value = "  +2.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01606"
$ws.Range("E38").Value = "  -0.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.914"
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8609"
$ws.Range("E40").Value = "  +1.67%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.011"
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.044.64"
$ws.Range("E42").Value = "  +2.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.27"
$ws.Range("E43").Value = "  +2.50%  "
$ws.Range("D44").Value = "1.805.28"
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.39"
$ws.Range("E45").Value = "  +0.36%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.017"
$ws.Range("E46").Value = "  +1.01%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈106"
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4385"
$ws.Range("E48").Value = "  +2.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.007"
$ws.Range("E49").Value = "  +2.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05165"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.440"
$ws.Range("E51").Value = "  -2.92%  "
